$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.522.37'
$ws.Range("E2").Value = '  -2.72%  '
$ws.Range("D3").Value = '2.422.11'
$ws.Range("E3").Value = '  -4.60%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.20'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.64'
$ws.Range("E6").Value = '  -5.41%  '
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.558'
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0968'
$ws.Range("E9").Value = '  -2.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.90'
$ws.Range("E11").Value = '  -6.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.318'
$ws.Range("E12").Value = '  -5.22%  '
$ws.Range("D13").Value = '2.869.05'
$ws.Range("E13").Value = '  -3.87%  '
$ws.Range("D14").Value = '57.499.16'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.39'
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("E16").Value = '  -4.19%  '
$ws.Range("D17").Value = '2.442.61'
$ws.Range("E17").Value = '  -3.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.31'
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.06'
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '309.46'
$ws.Range("E20").Value = '  -4.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.04'
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.51'
$ws.Range("E23").Value = '  -1.16%  '
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '2.559.29'
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.398'
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.155'
$ws.Range("E27").Value = '  -3.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.16'
$ws.Range("E28").Value = '  -4.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.06'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("D30").Value = '0.0₃0727'
$ws.Range("E30").Value = '  -4.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.67'
$ws.Range("E31").Value = '  -3.75%  '
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("E33").Value = '  -9.19%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.67'
$ws.Range("E36").Value = '  -3.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.17'
$ws.Range("E37").Value = '  -7.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.72'
$ws.Range("E38").Value = '  -6.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.14'
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.785'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  -6.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.35'
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.576'
$ws.Range("E43").Value = '  -4.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.70'
$ws.Range("E44").Value = '  -7.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0915'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '248.34'
$ws.Range("E46").Value = '  -11.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.85'
$ws.Range("E47").Value = '  -11.89%  '
$ws.Range("E48").Value = '  -3.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0208'
$ws.Range("E49").Value = '  -4.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.85'
$ws.Range("E50").Value = '  -5.71%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.680.31'
$ws.Range("E51").Value = '  -4.37%  '
